$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("B1").Value = "nodeid_linkid"
$ws.Range("B2").Select()
